$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 368, shifting existing rows 368:482 down to 369:483
$ws.Rows.Item(368).Insert()

# Populate the new row 368 with the new weekly data point
$ws.Range("A368").Value = 8
$ws.Range("B368").Value = "Terminal La Palmera de La Serena"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 45093
$ws.Range("E368").Value = 4
$ws.Range("F368").Value = 100112003
$ws.Range("G368").Value = "Ajo"
$ws.Range("H368").Value = "Chino"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 360
$ws.Range("K368").Value = 17500
$ws.Range("L368").Value = 18000
$ws.Range("M368").Value = 17750
$ws.Range("N368").Value = "`$/caja 10 kilos"
$ws.Range("O368").Value = "China"
$ws.Range("P368").Value = 1775
$ws.Range("Q368").Value = 10
$ws.Range("R368").Value = "Hortaliza"

Write-Output "done"
